$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 150 already has date (A150) and policy rate (G150); fill in B150:F150
$ws.Range("B150").Value = 145000
$ws.Range("C150").Value = 0.75
$ws.Range("D150").Value = 0.75
$ws.Range("E150").Value = 0.75
$ws.Range("F150").Value = 5

# New rows 151:172 -- date text (col A), monto/tasas/participantes (B:F), tasa politica (G)
# Use a throwaway Text-formatted cell to force Excel to keep the dd-mm-yyyy strings
# as literal text (shared string) instead of auto-converting them to date serials,
# then reset the style back to Normal so no extra style is applied to the cell.

$a = $ws.Range("A151")
$a.NumberFormat = "@"
$a.Value = "06-08-2021"
$a.Style = "Normal"
$ws.Range("B151").Value = 105000
$ws.Range("C151").Value = 0.75
$ws.Range("D151").Value = 0.75
$ws.Range("E151").Value = 0.75
$ws.Range("F151").Value = 4
$ws.Range("G151").Value = 0.75

$a = $ws.Range("A152")
$a.NumberFormat = "@"
$a.Value = "09-08-2021"
$a.Style = "Normal"
$ws.Range("B152").Value = 300000
$ws.Range("C152").Value = 0.75
$ws.Range("D152").Value = 0.75
$ws.Range("E152").Value = 0.75
$ws.Range("F152").Value = 4
$ws.Range("G152").Value = 0.75

$a = $ws.Range("A153")
$a.NumberFormat = "@"
$a.Value = "10-08-2021"
$a.Style = "Normal"
$ws.Range("B153").Value = 0
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0.75

$a = $ws.Range("A154")
$a.NumberFormat = "@"
$a.Value = "11-08-2021"
$a.Style = "Normal"
$ws.Range("B154").Value = 75000
$ws.Range("C154").Value = 0.75
$ws.Range("D154").Value = 0.75
$ws.Range("E154").Value = 0.75
$ws.Range("F154").Value = 4
$ws.Range("G154").Value = 0.75

$a = $ws.Range("A155")
$a.NumberFormat = "@"
$a.Value = "12-08-2021"
$a.Style = "Normal"
$ws.Range("B155").Value = 40000
$ws.Range("C155").Value = 0.75
$ws.Range("D155").Value = 0.75
$ws.Range("E155").Value = 0.75
$ws.Range("F155").Value = 3
$ws.Range("G155").Value = 0.75

$a = $ws.Range("A156")
$a.NumberFormat = "@"
$a.Value = "13-08-2021"
$a.Style = "Normal"
$ws.Range("B156").Value = 120000
$ws.Range("C156").Value = 0.75
$ws.Range("D156").Value = 0.75
$ws.Range("E156").Value = 0.75
$ws.Range("F156").Value = 5
$ws.Range("G156").Value = 0.75

$a = $ws.Range("A157")
$a.NumberFormat = "@"
$a.Value = "16-08-2021"
$a.Style = "Normal"
$ws.Range("B157").Value = 140000
$ws.Range("C157").Value = 0.75
$ws.Range("D157").Value = 0.75
$ws.Range("E157").Value = 0.75
$ws.Range("F157").Value = 6
$ws.Range("G157").Value = 0.75

$a = $ws.Range("A158")
$a.NumberFormat = "@"
$a.Value = "17-08-2021"
$a.Style = "Normal"
$ws.Range("B158").Value = 275000
$ws.Range("C158").Value = 0.75
$ws.Range("D158").Value = 0.75
$ws.Range("E158").Value = 0.75
$ws.Range("F158").Value = 8
$ws.Range("G158").Value = 0.75

$a = $ws.Range("A159")
$a.NumberFormat = "@"
$a.Value = "18-08-2021"
$a.Style = "Normal"
$ws.Range("B159").Value = 55000
$ws.Range("C159").Value = 0.75
$ws.Range("D159").Value = 0.75
$ws.Range("E159").Value = 0.75
$ws.Range("F159").Value = 3
$ws.Range("G159").Value = 0.75

$a = $ws.Range("A160")
$a.NumberFormat = "@"
$a.Value = "19-08-2021"
$a.Style = "Normal"
$ws.Range("B160").Value = 60000
$ws.Range("C160").Value = 0.75
$ws.Range("D160").Value = 0.75
$ws.Range("E160").Value = 0.75
$ws.Range("F160").Value = 3
$ws.Range("G160").Value = 0.75

$a = $ws.Range("A161")
$a.NumberFormat = "@"
$a.Value = "20-08-2021"
$a.Style = "Normal"
$ws.Range("B161").Value = 310000
$ws.Range("C161").Value = 0.75
$ws.Range("D161").Value = 0.75
$ws.Range("E161").Value = 0.75
$ws.Range("F161").Value = 7
$ws.Range("G161").Value = 0.75

$a = $ws.Range("A162")
$a.NumberFormat = "@"
$a.Value = "23-08-2021"
$a.Style = "Normal"
$ws.Range("B162").Value = 25000
$ws.Range("F162").Value = 2
$ws.Range("G162").Value = 0.75

$a = $ws.Range("A163")
$a.NumberFormat = "@"
$a.Value = "24-08-2021"
$a.Style = "Normal"
$ws.Range("B163").Value = 60000
$ws.Range("C163").Value = 0.75
$ws.Range("D163").Value = 0.75
$ws.Range("E163").Value = 0.75
$ws.Range("F163").Value = 3
$ws.Range("G163").Value = 0.75

$a = $ws.Range("A164")
$a.NumberFormat = "@"
$a.Value = "25-08-2021"
$a.Style = "Normal"
$ws.Range("B164").Value = 130000
$ws.Range("C164").Value = 0.75
$ws.Range("D164").Value = 0.75
$ws.Range("E164").Value = 0.75
$ws.Range("F164").Value = 4
$ws.Range("G164").Value = 0.75

$a = $ws.Range("A165")
$a.NumberFormat = "@"
$a.Value = "26-08-2021"
$a.Style = "Normal"
$ws.Range("B165").Value = 160000
$ws.Range("C165").Value = 0.75
$ws.Range("D165").Value = 0.75
$ws.Range("E165").Value = 0.75
$ws.Range("F165").Value = 5
$ws.Range("G165").Value = 0.75

$a = $ws.Range("A166")
$a.NumberFormat = "@"
$a.Value = "27-08-2021"
$a.Style = "Normal"
$ws.Range("B166").Value = 130000
$ws.Range("C166").Value = 0.75
$ws.Range("D166").Value = 0.75
$ws.Range("E166").Value = 0.75
$ws.Range("F166").Value = 5
$ws.Range("G166").Value = 0.75

$a = $ws.Range("A167")
$a.NumberFormat = "@"
$a.Value = "30-08-2021"
$a.Style = "Normal"
$ws.Range("B167").Value = 525000
$ws.Range("C167").Value = 0.75
$ws.Range("D167").Value = 0.75
$ws.Range("E167").Value = 0.75
$ws.Range("F167").Value = 8
$ws.Range("G167").Value = 0.75

$a = $ws.Range("A168")
$a.NumberFormat = "@"
$a.Value = "31-08-2021"
$a.Style = "Normal"
$ws.Range("B168").Value = 0
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0.75

$a = $ws.Range("A169")
$a.NumberFormat = "@"
$a.Value = "01-09-2021"
$a.Style = "Normal"
$ws.Range("B169").Value = 0
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 1.5

$a = $ws.Range("A170")
$a.NumberFormat = "@"
$a.Value = "02-09-2021"
$a.Style = "Normal"
$ws.Range("B170").Value = 0
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 1.5

$a = $ws.Range("A171")
$a.NumberFormat = "@"
$a.Value = "03-09-2021"
$a.Style = "Normal"
$ws.Range("B171").Value = 0
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 1.5

$a = $ws.Range("A172")
$a.NumberFormat = "@"
$a.Value = "06-09-2021"
$a.Style = "Normal"
$ws.Range("G172").Value = 1.5
